$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.724581718444824
$ws.Range("B1").Value = 4.606094360351562
$ws.Range("C1").Value = 3.125218152999878
$ws.Range("D1").Value = 1.968222856521606
$ws.Range("E1").Value = 1.455422759056091
